# The document's single table holds 20 rows x 5 columns of
# "a+b=" / "a-b=" practice problems. The commit replaces the
# expression text in every cell (some old expressions repeat, so we
# must target each cell by position rather than doing a single global
# Find/Replace). $pairs lists (oldText, newText) in row-major cell
# order matching the table layout.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$pairs = @(
    @("64-56=", "6+56="),
    @("49+29=", "18-13="),
    @("0+29=", "38+22="),
    @("49+23=", "71-29="),
    @("47-44=", "85-62="),
    @("28-8=", "83-7="),
    @("14+79=", "74-0="),
    @("75-48=", "28+29="),
    @("75-3=", "18+42="),
    @("30-11=", "6+82="),
    @("4+86=", "72-45="),
    @("47+21=", "12+46="),
    @("93-49=", "50-21="),
    @("58-2=", "15+38="),
    @("16+8=", "82-60="),
    @("88-66=", "57+30="),
    @("40+3=", "57-12="),
    @("57+9=", "37-22="),
    @("5+89=", "46-22="),
    @("80-1=", "26-10="),
    @("89+7=", "43-43="),
    @("77+9=", "49-25="),
    @("83-49=", "64+5="),
    @("18-18=", "49-8="),
    @("75-51=", "72-13="),
    @("92-43=", "80-23="),
    @("68-35=", "5+56="),
    @("57-46=", "94-60="),
    @("15+55=", "23+14="),
    @("23-14=", "21-6="),
    @("58-16=", "61+19="),
    @("54-7=", "91-8="),
    @("27+20=", "65+5="),
    @("94-64=", "29-26="),
    @("74-23=", "24+13="),
    @("65-15=", "20+42="),
    @("46-7=", "87-40="),
    @("98-27=", "84-13="),
    @("68+24=", "92-29="),
    @("44+41=", "65-10="),
    @("90-82=", "99-27="),
    @("55-4=", "4+81="),
    @("60+19=", "23+20="),
    @("19+52=", "1+91="),
    @("2+81=", "61+11="),
    @("51-35=", "32+0="),
    @("65+13=", "63-34="),
    @("68+13=", "39+48="),
    @("86-46=", "7+64="),
    @("86-69=", "29+52="),
    @("71-52=", "31-30="),
    @("67+6=", "12+82="),
    @("24-3=", "92+5="),
    @("13+64=", "54-34="),
    @("11+88=", "43+23="),
    @("48-7=", "49+18="),
    @("57-44=", "0+95="),
    @("11+60=", "89-56="),
    @("65+19=", "37-13="),
    @("85-27=", "79+17="),
    @("26-16=", "85-59="),
    @("72-20=", "78-60="),
    @("99-29=", "14+76="),
    @("49-38=", "51-34="),
    @("28+37=", "38+32="),
    @("76-6=", "1+58="),
    @("31+46=", "77+3="),
    @("49-9=", "99-74="),
    @("57+33=", "64+27="),
    @("97-72=", "99-57="),
    @("28+1=", "70-70="),
    @("19+12=", "39+13="),
    @("30-25=", "3+55="),
    @("80-40=", "72-61="),
    @("34-33=", "18+49="),
    @("12+59=", "70-10="),
    @("94-28=", "18-17="),
    @("66-15=", "86-13="),
    @("37-8=", "12+8="),
    @("93-42=", "44-32="),
    @("47+33=", "29+49="),
    @("6+30=", "22+22="),
    @("84+11=", "10+44="),
    @("1+42=", "3+73="),
    @("86-33=", "61+26="),
    @("33-22=", "33+34="),
    @("16+20=", "95-25="),
    @("54-10=", "56-16="),
    @("39-11=", "39+5="),
    @("68+24=", "12-3="),
    @("46+16=", "99-98="),
    @("40-24=", "91-59="),
    @("3+81=", "62-28="),
    @("30+51=", "92-60="),
    @("40-12=", "34+40="),
    @("5+38=", "3+66="),
    @("1+44=", "99-31="),
    @("36+54=", "75-59="),
    @("65-47=", "87-76="),
    @("15+50=", "53+14=")
)

$rows = 20
$cols = 5
$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $pair = $pairs[$i]
        $old = $pair[0]
        $new = $pair[1]

        # Build a plain Document.Range from the cell's Start/End instead of
        # using Cell.Range directly: Cell.Range includes the end-of-cell
        # marker and, more importantly, Find.Execute on it isn't reliably
        # confined to the cell - it can touch matching text elsewhere in the
        # document. A Document.Range(start, end) plus wdReplaceOne (1) keeps
        # the replacement pinned to this exact cell, which matters because
        # several "old" expressions (e.g. "68+24=") occur more than once.
        $cellRange = $t.Cell($r, $c).Range
        $rng = $d.Range($cellRange.Start, $cellRange.End)
        $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
        if (-not $found) {
            Write-Host "MISS row=$r col=$c old=$old new=$new"
        }
        $i = $i + 1
    }
}
Write-Host "Done, processed $i cells"
